$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Handback report generation: the zh-cn and de-de detail sheets gain two
# new columns of data (F = Latest Target File, G = Latest Handback File)
# for each of the two tracked source files, the Status column moves from
# "Ready for handoff" to "Handed back: in sync with en-US", and the
# Latest Handback DateTime column (H) is stamped with the real handback
# time instead of the zero-date placeholder.
# ---------------------------------------------------------------------

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet -----------------------------------------------------
# The per-language status columns summarise the same handoff status text
# used on the detail sheets, so they pick up the new wording too.

$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet ------------------------------------------------------

$wsZh.Range("C2").Value = $statusHandedBack
$wsZh.Range("C3").Value = $statusHandedBack

$wsZh.Range("F2").Value = "39b07019-896a-4d16-842b-bb42829f0703.md"
$wsZh.Range("G2").Value = "39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.zh-cn.xlf"

$wsZh.Range("F3").Value = "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md"
$wsZh.Range("G3").Value = "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.zh-cn.xlf"

$wsZh.Range("H2").Value = "2016-03-13 06:50:14"
$wsZh.Range("H3").Value = "2016-03-13 06:50:14"

# Rebuild the hyperlinks in display order so F/G land alongside A/B/D.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/39b07019-896a-4d16-842b-bb42829f0703.md", "", "", "39b07019-896a-4d16-842b-bb42829f0703.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/39b07019-896a-4d16-842b-bb42829f0703.md", "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fa66c6db28e2bbb0ef480c9bfc26f1bdc51bb086/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.zh-cn.xlf", "", "", "39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/39b07019-896a-4d16-842b-bb42829f0703.md", "", "", "39b07019-896a-4d16-842b-bb42829f0703.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fa66c6db28e2bbb0ef480c9bfc26f1bdc51bb086/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.zh-cn.xlf", "", "", "39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.zh-cn.xlf") | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md", "", "", "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md", "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fa66c6db28e2bbb0ef480c9bfc26f1bdc51bb086/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.zh-cn.xlf", "", "", "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md", "", "", "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fa66c6db28e2bbb0ef480c9bfc26f1bdc51bb086/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.zh-cn.xlf", "", "", "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.zh-cn.xlf") | Out-Null

foreach ($addr in @("A2","B2","D2","F2","G2","A3","B3","D3","F3","G3")) {
    $wsZh.Range($addr).Font.Underline = 2
    $wsZh.Range($addr).Font.Color = 15570276
}

# --- de-de sheet --------------------------------------------------------

$wsDe.Range("C2").Value = $statusHandedBack
$wsDe.Range("C3").Value = $statusHandedBack

$wsDe.Range("F2").Value = "39b07019-896a-4d16-842b-bb42829f0703.md"
$wsDe.Range("G2").Value = "39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.de-de.xlf"

$wsDe.Range("F3").Value = "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md"
$wsDe.Range("G3").Value = "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.de-de.xlf"

$wsDe.Range("H2").Value = "2016-03-13 06:50:21"
$wsDe.Range("H3").Value = "2016-03-13 06:50:21"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/39b07019-896a-4d16-842b-bb42829f0703.md", "", "", "39b07019-896a-4d16-842b-bb42829f0703.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/39b07019-896a-4d16-842b-bb42829f0703.md", "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01aa9f03f0dbd5ab36c898cd118e5fc6c73d2b01/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.de-de.xlf", "", "", "39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/39b07019-896a-4d16-842b-bb42829f0703.md", "", "", "39b07019-896a-4d16-842b-bb42829f0703.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01aa9f03f0dbd5ab36c898cd118e5fc6c73d2b01/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.de-de.xlf", "", "", "39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.de-de.xlf") | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md", "", "", "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md", "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01aa9f03f0dbd5ab36c898cd118e5fc6c73d2b01/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.de-de.xlf", "", "", "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md", "", "", "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01aa9f03f0dbd5ab36c898cd118e5fc6c73d2b01/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.de-de.xlf", "", "", "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.de-de.xlf") | Out-Null

foreach ($addr in @("A2","B2","D2","F2","G2","A3","B3","D3","F3","G3")) {
    $wsDe.Range($addr).Font.Underline = 2
    $wsDe.Range($addr).Font.Color = 15570276
}
